# "add runner of deploy"
#
# The test-case sheet has a "执行者" (runner/executor) column (H) that was
# left blank for every test row. This fills it in for all data rows
# (2-43) with the runner's name, picking up the same visual format/
# validation already used on the neighbouring "开发者" (developer) column
# F (F2's format: centered, bottom-aligned, wrapped text, same border/
# fill), and extends the data-validation rule to the newly filled cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Borrow F2's cell formatting (font/fill/border/alignment) for the whole
# runner column so the new entries look consistent with the rest of the
# table.
$ws.Range("F2").Copy()
$ws.Range("H2:H43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Record who ran the tests.
$ws.Range("H2:H43").Value = "刘彩丽"

# Column F carries a (pass-through) data validation rule; mirror it onto
# the runner column now that it holds real data too.
$ws.Range("H2:H43").Validation.Add(8, 1, 1, "0", "0")

# Leave the selection where the edit finished, like a user tabbing/filling
# down to the last row.
$ws.Range("H42:H43").Select()
